$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(22)
$ws.Copy($null, $ws)
$wsNew1 = $wb.Worksheets.Item(23)
$wsNew1.Name = "VT-AuthCapVoid-Generic"

$wsNew1.Copy($null, $wsNew1)
$wsNew2 = $wb.Worksheets.Item(24)
$wsNew2.Name = "VT-AuthCapCredit-Generic"

$wsNew2.Rows.Item(4).Resize(2).Delete()
$wsNew2.Columns.AutoFit()
Write-Output $wsNew2.Columns.Item(5).ColumnWidth
Write-Output $wsNew2.Columns.Item(6).ColumnWidth
